# Update the "K" column (column G) with newly regenerated strikeout (K)
# values, replacing the old "Strike#" derived values, per the save_data
# regeneration described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 2
    6  = 3
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 1
    16 = 2
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 2
    22 = 2
    23 = 3
    24 = 1
    25 = 0
    26 = 1
    27 = 0
    28 = 2
    29 = 2
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
